$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to reset style
# after forcing a text number-format on price cells so no stray style is introduced.
$defaultStyle = $ws.Range("D4").Style

# Row 2: Price
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.849.58'
$ws.Range("D2").Style = $defaultStyle
# Row 2: Volume(1h)
$ws.Range("E2").Value = '  +0.78%  '

# Row 3: Price
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.140.39'
$ws.Range("D3").Style = $defaultStyle
# Row 3: Volume(1h)
$ws.Range("E3").Value = '  +1.66%  '

# Row 5: Price
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.06'
$ws.Range("D5").Style = $defaultStyle
# Row 5: Volume(1h)
$ws.Range("E5").Value = '  +1.85%  '

# Row 6: Price
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.77'
$ws.Range("D6").Style = $defaultStyle
# Row 6: Volume(1h)
$ws.Range("E6").Value = '  +3.88%  '

# Row 7: Volume(1h)
$ws.Range("E7").Value = '  -0.07%  '

# Row 8: Price
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.136.65'
$ws.Range("D8").Style = $defaultStyle
# Row 8: Volume(1h)
$ws.Range("E8").Value = '  +1.64%  '

# Row 9: Volume(1h)
$ws.Range("E9").Value = '  +3.95%  '

# Row 10: Volume(1h)
$ws.Range("E10").Value = '  +4.57%  '

# Row 11: Volume(1h)
$ws.Range("E11").Value = '  +0.56%  '

# Row 12: Volume(1h)
$ws.Range("E12").Value = '  +6.18%  '

# Row 13: Volume(1h)
$ws.Range("E13").Value = '  +10.72%  '

# Row 14: Price
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.39'
$ws.Range("D14").Style = $defaultStyle
# Row 14: Volume(1h)
$ws.Range("E14").Value = '  +6.17%  '

# Row 15: Price
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.653.32'
$ws.Range("D15").Style = $defaultStyle
# Row 15: Volume(1h)
$ws.Range("E15").Value = '  +1.93%  '

# Row 16: Price
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.913.25'
$ws.Range("D16").Style = $defaultStyle
# Row 16: Volume(1h)
$ws.Range("E16").Value = '  +0.82%  '

# Row 17: Volume(1h)
$ws.Range("E17").Value = '  +6.30%  '

# Row 18: Price
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.139.87'
$ws.Range("D18").Style = $defaultStyle
# Row 18: Volume(1h)
$ws.Range("E18").Value = '  +1.73%  '

# Row 19: Volume(1h)
$ws.Range("E19").Value = '  +0.28%  '

# Row 20: Price
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '511.70'
$ws.Range("D20").Style = $defaultStyle
# Row 20: Volume(1h)
$ws.Range("E20").Value = '  +6.51%  '

# Row 21: Price
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.91'
$ws.Range("D21").Style = $defaultStyle
# Row 21: Volume(1h)
$ws.Range("E21").Value = '  +7.03%  '

# Row 22: Volume(1h)
$ws.Range("E22").Value = '  +8.74%  '

# Row 23: Price
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.30'
$ws.Range("D23").Style = $defaultStyle
# Row 23: Volume(1h)
$ws.Range("E23").Value = '  +10.48%  '

# Row 24: Volume(1h)
$ws.Range("E24").Value = '  +3.50%  '

# Row 25: Price
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.90'
$ws.Range("D25").Style = $defaultStyle
# Row 25: Volume(1h)
$ws.Range("E25").Value = '  +4.50%  '

# Row 26: Volume(1h)
$ws.Range("E26").Value = '  +0.46%  '

# Row 27: Volume(1h)
$ws.Range("E27").Value = '  +3.38%  '

# Row 28: Volume(1h)
$ws.Range("E28").Value = '  +8.23%  '

# Row 29: Volume(1h)
$ws.Range("E29").Value = '  +4.98%  '

# Row 30: Price
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.91'
$ws.Range("D30").Style = $defaultStyle
# Row 30: Volume(1h)
$ws.Range("E30").Value = '  +6.27%  '

# Row 31: Price
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = $defaultStyle
# Row 31: Volume(1h)
$ws.Range("E31").Value = '  +0.04%  '

# Row 32: Volume(1h)
$ws.Range("E32").Value = '  +3.43%  '

# Row 33: Price
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.65'
$ws.Range("D33").Style = $defaultStyle
# Row 33: Volume(1h)
$ws.Range("E33").Value = '  +6.03%  '

# Row 34: Price
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.06'
$ws.Range("D34").Style = $defaultStyle
# Row 34: Volume(1h)
$ws.Range("E34").Value = '  +7.93%  '

# Row 35: Volume(1h)
$ws.Range("E35").Value = '  +5.90%  '

# Row 36: Price
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.39'
$ws.Range("D36").Style = $defaultStyle
# Row 36: Volume(1h)
$ws.Range("E36").Value = '  -0.76%  '

# Row 37: Price
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '478.94'
$ws.Range("D37").Style = $defaultStyle
# Row 37: Volume(1h)
$ws.Range("E37").Value = '  +4.28%  '

# Row 38: Volume(1h)
$ws.Range("E38").Value = '  +3.50%  '

# Row 39: Volume(1h)
$ws.Range("E39").Value = '  +4.38%  '

# Row 40: Price
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("D40").Style = $defaultStyle
# Row 40: Volume(1h)
$ws.Range("E40").Value = '  -0.90%  '

# Row 41: Price
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.113.13'
$ws.Range("D41").Style = $defaultStyle
# Row 41: Volume(1h)
$ws.Range("E41").Value = '  +4.22%  '

# Row 42: Price
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.63'
$ws.Range("D42").Style = $defaultStyle
# Row 42: Volume(1h)
$ws.Range("E42").Value = '  +4.35%  '

# Row 43: Volume(1h)
$ws.Range("E43").Value = '  +4.05%  '

# Row 44: Volume(1h)
$ws.Range("E44").Value = '  +11.41%  '

# Row 45: Price
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").Style = $defaultStyle
# Row 45: Volume(1h)
$ws.Range("E45").Value = '  +14.34%  '

# Row 46: Price
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.10'
$ws.Range("D46").Style = $defaultStyle
# Row 46: Volume(1h)
$ws.Range("E46").Value = '  +3.63%  '

# Row 47: Price
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0570'
$ws.Range("D47").Style = $defaultStyle
# Row 47: Volume(1h)
$ws.Range("E47").Value = '  +10.15%  '

# Row 49: Price
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.116'
$ws.Range("D49").Style = $defaultStyle
# Row 49: Volume(1h)
$ws.Range("E49").Value = '  +3.19%  '

# Row 50: Volume(1h)
$ws.Range("E50").Value = '  +10.03%  '

# Row 51: Price
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.74'
$ws.Range("D51").Style = $defaultStyle

